$d = $word.ActiveDocument

$replacements = @(
    @("891÷3=", "170÷4="),
    @("983÷4=", "381÷3="),
    @("637÷2=", "602÷2="),
    @("948÷3=", "128÷3="),
    @("426÷9=", "339÷8="),
    @("834÷2=", "559÷5="),
    @("227÷9=", "916÷2="),
    @("332÷9=", "268÷7="),
    @("678÷6=", "718÷2="),
    @("589÷3=", "305÷3="),
    @("197÷8=", "320÷8="),
    @("148÷9=", "115÷4="),
    @("249÷5=", "314÷8="),
    @("169÷6=", "715÷7="),
    @("377÷6=", "771÷6="),
    @("786÷5=", "563÷4="),
    @("731÷5=", "338÷4="),
    @("974÷9=", "924÷4="),
    @("359÷6=", "343÷4="),
    @("669÷7=", "860÷6="),
    @("459÷5=", "362÷8="),
    @("698÷7=", "189÷2="),
    @("919÷3=", "308÷2="),
    @("517÷3=", "956÷6="),
    @("405÷3=", "520÷8="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
